# Apply change: insert a new "Industry" column (C) between "Stock Name" (B)
# and "Mutual Fund" (previously C, now D), shifting existing columns C:I to D:J.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C; this shifts C..I -> D..J and copies
# formatting from the column to the right (C), matching the original style.
$ws.Columns.Item(3).Insert()

# Set header for the new Industry column, matching the style used by the
# other header cells in row 1 (bold/centered/bordered).
$ws.Cells.Item(1, 3).Value = "Industry"
$ws.Cells.Item(1, 4).Copy()
$ws.Cells.Item(1, 3).PasteSpecial(-4122)
$ws.Cells.Item(1, 3).Value = "Industry"

# Populate the Industry values for each data row (2-32).
$industries = @{
    2  = "Banks"
    3  = "Finance"
    4  = "Auto Components"
    5  = "Power"
    6  = "Pharmaceuticals & Biotechnology"
    7  = "Insurance"
    8  = "Banks"
    9  = "Food Products"
    10 = "Metals & Minerals Trading"
    11 = "Personal Products"
    12 = "Capital Markets"
    13 = "Pharmaceuticals & Biotechnology"
    14 = "Cement & Cement Products"
    15 = "Personal Products"
    16 = "Consumable Fuels"
    17 = "Finance"
    18 = "Entertainment"
    19 = "Telecom - Services"
    20 = "IT - Software"
    21 = "Automobiles"
    22 = "Realty"
    23 = "Pharmaceuticals & Biotechnology"
    24 = "Power"
    25 = "Petroleum Products"
    26 = "Textiles & Apparels"
    27 = "Healthcare"
    28 = "Retailing"
    29 = "Insurance"
    30 = "Pharmaceuticals & Biotechnology"
    31 = "Petroleum Products"
    32 = "Pharmaceuticals & Biotechnology"
}

foreach ($row in $industries.Keys) {
    $ws.Cells.Item($row, 3).Value = $industries[$row]
}
